# Delete row 44 from the "AYKO" sheet, shifting all rows below it up by one.
# This corresponds to removing the "-337 / PARAGUAY /ALT/ 4259" record and
# shrinking the used range from A1:P112 down to A1:P111.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("AYKO")

$ws.Rows.Item(44).Delete()
